$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to text so they
# keep matching the original inline-string ("Price" column) representation.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "26.280.90"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.607.13"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "212.88"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.487"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "18.50"
$ws.Range("E10").Value = "  +2.61%  "
$ws.Range("D11").Value = "0.0813"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "1.830.92"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "1.605.24"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "26.260.83"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "62.11"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "200.87"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").Value = "143.72"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("E30").Value = "  +4.98%  "
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").Value = "3.21"
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("D35").Value = "2.38"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("D36").Value = "1.163.55"
$ws.Range("E36").Value = "  +3.60%  "
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").Value = "0.788"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("E41").Value = "  +0.89%  "
$ws.Range("D42").Value = "5.36"
$ws.Range("E42").Value = "  +4.37%  "
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("D44").Value = "1.742.05"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").Value = "92.08"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").Value = "0.0₆0105"
$ws.Range("E47").Value = "  +13.84%  "
$ws.Range("D48").Value = "54.08"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("E51").Value = "  -0.08%  "
